$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Wins, Losses, Ties in AD1:AF1, matching style of existing headers (e.g. AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-43: Wins=93, Losses=69, Ties=0
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 93
    $ws.Cells.Item($r, 31).Value = 69
    $ws.Cells.Item($r, 32).Value = 0
}
